$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be read as numbers
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = '79.750.12'
$ws.Range("E2").Value = '  +4.68%  '
$ws.Range("D3").Value = '3.167.56'
$ws.Range("E3").Value = '  +3.14%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '207.82'
$ws.Range("E5").Value = '  +4.97%  '
$ws.Range("D6").Value = '628.66'
$ws.Range("E6").Value = '  +2.07%  '
$ws.Range("D7").Value = '0.270'
$ws.Range("E7").Value = '  +29.81%  '
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("E9").Value = '  +8.30%  '
$ws.Range("D10").Value = '3.168.40'
$ws.Range("E10").Value = '  +3.21%  '
$ws.Range("D11").Value = '0.608'
$ws.Range("E11").Value = '  +38.36%  '
$ws.Range("D12").Value = '0.0000256'
$ws.Range("E12").Value = '  +32.90%  '
$ws.Range("E13").Value = '  +2.68%  '
$ws.Range("B14").Value = 'Toncoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D14").Value = '5.29'
$ws.Range("E14").Value = '  +1.51%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '3.760.33'
$ws.Range("E15").Value = '  +4.17%  '
$ws.Range("D16").Value = '32.00'
$ws.Range("E16").Value = '  +10.56%  '
$ws.Range("D17").Value = '79.834.60'
$ws.Range("E17").Value = '  +4.86%  '
$ws.Range("D18").Value = '3.192.70'
$ws.Range("E18").Value = '  +3.88%  '
$ws.Range("D19").Value = '14.39'
$ws.Range("E19").Value = '  +6.17%  '
$ws.Range("D20").Value = '442.54'
$ws.Range("E20").Value = '  +16.66%  '
$ws.Range("B21").Value = 'SuiNetwork'
$ws.Range("C21").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D21").Value = '2.95'
$ws.Range("E21").Value = '  +18.89%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '9.16'
$ws.Range("E22").Value = '  +0.71%  '
$ws.Range("D23").Value = '5.32'
$ws.Range("E23").Value = '  +20.83%  '
$ws.Range("B24").Value = 'LEO'
$ws.Range("C24").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D24").Value = '6.90'
$ws.Range("E24").Value = '  +6.63%  '
$ws.Range("B25").Value = 'WrappedeETH'
$ws.Range("C25").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D25").Value = '3.351.21'
$ws.Range("E25").Value = '  +3.68%  '
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").Value = '76.78'
$ws.Range("E26").Value = '  +6.53%  '
$ws.Range("B27").Value = 'NEARProtocol'
$ws.Range("C27").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D27").Value = '4.72'
$ws.Range("E27").Value = '  +9.11%  '
$ws.Range("B28").Value = 'Aptos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D28").Value = '10.91'
$ws.Range("E28").Value = '  +10.72%  '
$ws.Range("B29").Value = 'Dai'
$ws.Range("C29").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.16%  '
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0000123'
$ws.Range("E30").Value = '  +13.77%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '9.14'
$ws.Range("E31").Value = '  +10.20%  '
$ws.Range("B32").Value = 'Binance-PegBSC-USD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  +0.53%  '
$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").Value = '549.58'
$ws.Range("E33").Value = '  +10.03%  '
$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").Value = '1.48'
$ws.Range("E34").Value = '  +4.92%  '
$ws.Range("B35").Value = 'PancakeSwap'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D35").Value = '2.00'
$ws.Range("E35").Value = '  +4.87%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = '0.150'
$ws.Range("E36").Value = '  +21.69%  '
$ws.Range("B37").Value = 'EthereumClassic'
$ws.Range("C37").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D37").Value = '23.34'
$ws.Range("E37").Value = '  +12.73%  '
$ws.Range("B38").Value = 'Cronos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D38").Value = '0.122'
$ws.Range("E38").Value = '  +19.14%  '
$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  -0.04%  '
$ws.Range("B40").Value = 'PolygonEcosystemToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D40").Value = '0.409'
$ws.Range("E40").Value = '  +8.28%  '
$ws.Range("B41").Value = 'WhiteBITCoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D41").Value = '20.78'
$ws.Range("E41").Value = '  +3.61%  '
$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D42").Value = '164.74'
$ws.Range("E42").Value = '  +1.25%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D43").Value = '5.67'
$ws.Range("E43").Value = '  +11.18%  '
$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '189.21'
$ws.Range("E45").Value = '  -2.99%  '
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").Value = '1.82'
$ws.Range("E46").Value = '  +11.07%  '
$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D47").Value = '2.70'
$ws.Range("E47").Value = '  +11.39%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").Value = '0.786'
$ws.Range("E48").Value = '  -1.34%  '
$ws.Range("D49").Value = '1.32'
$ws.Range("E49").Value = '  +5.85%  '
$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D50").Value = '43.53'
$ws.Range("E50").Value = '  +5.28%  '
$ws.Range("B51").Value = 'Filecoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D51").Value = '4.28'
$ws.Range("E51").Value = '  +10.37%  '
